$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column Y (12-jul), one day after the existing last column X (11-jul)
$ws.Range("Y1").Value = "12-jul"
$ws.Range("Y1").NumberFormat = $ws.Range("X1").NumberFormat

# Add the new day's values for each product row, mirroring style/format of column X
$values = @(13, 19, 6, 9, 14, 16, 16, 10, 20, 26)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $cell = $ws.Cells.Item($row, 25)  # column 25 = Y
    $cell.Value = $values[$i]
    $cell.HorizontalAlignment = $ws.Cells.Item($row, 24).HorizontalAlignment
    $cell.NumberFormat = $ws.Cells.Item($row, 24).NumberFormat
}

# Move the active selection to Y12 (one row below the new last data row), matching X12->X12->Y12 pattern
$ws.Range("Y12").Select()
